$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.470.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.13%  "
$ws.Range("D3").Value = "'2.997.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.79%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'543.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "'130.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.58%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'2.991.53"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.83%  "
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("D10").Value = "'5.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.39%  "
$ws.Range("D11").Value = "'0.144"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.84%  "
$ws.Range("E12").Value = "  -3.44%  "
$ws.Range("E13").Value = "  -3.27%  "
$ws.Range("D14").Value = "'33.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.43%  "
$ws.Range("D15").Value = "'3.484.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.83%  "
$ws.Range("D16").Value = "'61.592.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.98%  "
$ws.Range("D17").Value = "'0.109"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.94%  "
$ws.Range("D18").Value = "'3.001.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.65%  "
$ws.Range("D19").Value = "'6.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("D20").Value = "'478.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "'13.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.97%  "
$ws.Range("D22").Value = "'0.661"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.95%  "
$ws.Range("D23").Value = "'6.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.14%  "
$ws.Range("E24").Value = "  +2.63%  "
$ws.Range("D25").Value = "'11.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'2.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.67%  "
$ws.Range("D28").Value = "'7.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "'1.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").Value = "'25.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.10%  "
$ws.Range("D32").Value = "'1.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.90%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "'5.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").Value = "'54.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.93%  "
$ws.Range("D36").Value = "'5.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.49%  "
$ws.Range("D37").Value = "'442.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.52%  "
$ws.Range("D38").Value = "'3.124.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.15%  "
$ws.Range("D39").Value = "'0.0788"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").Value = "'0.0380"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.62%  "
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("D42").Value = "'8.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("D44").Value = "'2.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.59%  "
$ws.Range("D45").Value = "'25.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("E46").Value = "  -5.02%  "
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("E48").Value = "  -4.87%  "
$ws.Range("E49").Value = "  +9.00%  "
$ws.Range("D50").Value = "'113.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.75%  "
$ws.Range("D51").Value = "'0.0₃0481"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.44%  "
